# Auto-generated data application for rows 2..40 (A..I)
# productos.xlsx: refresh Dahua brand naming, fix row 12 placeholder data,
# and append new JAGL products (rows 13-40).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(1.0, "Domo Plástico *CON AUDIO* de 2MP, 20mts IR", 25.51, 198.978, 5.0, "CAMARA", "DAHUA", "", "📷"),
  @(2.0, "Cámara PTZ HDCVI Starlight de 2MP y 15x IK10 ", 370.51, 2890.0, 63.0, "CAMARA", "DAHUA", "", "📷"),
  @(3.0, "Bala Imou de 3MP (2304 x 1296), Lente de 2.8mm, IR hasta 30m ", 78.72, 614.0, 5.0, "CAMARA", "IMOU", "", "📷"),
  @(4.0, "Cámara 4MP de acceso y lectora de placas ANPR, clasificación IP67 e IK10 Dahua Varifocal ", 629.49, 4910.0, 3.0, "CAMARA", "DAHUA", "", "📷"),
  @(5.0, "IMOU GO KIT: 1x IPC-B32P-V2-IMOU imou Cámara Cell Go WiFi , 3MP, IR de 7m ,PIR, Micrófono, Altavoces, Batería Recargable y panel solar ", 128.07692307692307, 999.0, 1.0, "CAMARA", "IMOU", "", "📷"),
  @(6.0, "Imou Bombillo cámara 5MP WiFi para hogares ", 58.84615384615385, 459.0, 4.0, "CAMARA", "IMOU", "", "📷"),
  @(7.0, "DISCO DURO DE 1TB, *SEAGATE SKYHAWK* ", 660.0, 5148.0, 6.0, "SISTEMA DE ALMACENAMIENTO", "DAHUA", "", "💾"),
  @(8.0, "DISCO DURO DE 1TERABYTE B-STOCK ", 57.69230769230769, 450.0, 3.0, "SISTEMA DE ALMACENAMIENTO", "TOSHIBA", "", "💾"),
  @(9.0, "DISCO DURO DE 4TERABYTES B-STOCK ", 142.94871794871796, 1115.0, 7.0, "SISTEMA DE ALMACENAMIENTO", "SEAGATE", "", "💾"),
  @(10.0, "DISCO DURO DE 8TB, *SEAGATE SKYHAWK* ", 332.05128205128204, 2590.0, 6.0, "SISTEMA DE ALMACENAMIENTO", "SEAGATE ", "", "💾"),
  @(12.0, "Disco Duro de 10TB, *SEAGATE SKYHAWK DAHUA* ", 482.43589743589746, 3763.0, 3.0, "SISTEMA DE ALMACENAMIENTO", "SEAGATE ", "", "💾"),
  @(13.0, "Disco duro de estado sólido SATA de 2,5 pulgadas para grabado 24/7 de 1 TB ", 80.77, 630.01, 2.0, "SISTEMA DE ALMACENAMIENTO", "DAHUA", "", "💾"),
  @(14.0, "Disco duro de estado sólido SATA de 2,5 pulgadas para grabado 24/7 de 512 GB ", 56.92, 443.98, 5.0, "SISTEMA DE ALMACENAMIENTO", "DAHUA", "", "💾"),
  @(15.0, "NVR De 8 MP/4k de 32 Canales IP, 160 Mbps con 2 Bahías de discos duros de hasta 10 TB ", 287.82, 2245.0, 4.0, "SISTEMA DE ALMACENAMIENTO", "DAHUA", "", "💾"),
  @(16.0, "Bobina de 100 metros de cable UTP Cat 5e CCA, en color blanco anti flama ", 19.23, 150.0, 51.0, "ACCESORIO CCTV", "DAHUA", "", "🔌"),
  @(17.0, "Bobina de Cable UTP Cat6 100% Cobre Anti Flama ", 141.02564102564102, 1100.0, 4.0, "ACCESORIO CCTV", "DAHUA", "", "🔌"),
  @(18.0, "Bobina de cable *Coaxial RG59*, Sprywire ", 63.97, 499.0, 5.0, "ACCESORIO CCTV", "SPRYWIRE", "", "🔌"),
  @(19.0, "Bobina de cable *Exterior con Gel CAT5E* (305mts), SpryWire ", 83.97, 654.966, 8.0, "ACCESORIO CCTV", "SPRYWIRE", "", "🔌"),
  @(20.0, "Bobina de cable *Interior CAT6* (305mts), SpryWire ", 67.82051282051282, 529.0, 9.0, "ACCESORIO CCTV", "SPRYWIRE", "", "🔌"),
  @(21.0, "CABLE HDMI DE 15mts 4k 1.4V ", 25.512820512820515, 199.0, 8.0, "ACCESORIO CCTV", "SPRYWIRE", "", "🔌"),
  @(22.0, "CABLE HDMI DE 40mts 1080P ", 91.03, 710.0, 5.0, "ACCESORIO CCTV", "SPRYWIRE ", "", "🔌"),
  @(23.0, "AJAX HUB2 4G PANEL DE CONTROL ETHERNET CON SOPORTE PARA LA VERIFICACIÓN FOTOGRÁFICA DE ALARMAS NEGRO ", 384.62, 3000.0, 5.0, "CONTROL DE ACCESO Y SEGURIDAD", "AJAX ", "", "🔒"),
  @(24.0, "AJAX Hub: Panel de control de seguridad con soporte para la verificación fotográfica de alarmas ", 1802.0, 14055.6, 9.0, "CONTROL DE ACCESO Y SEGURIDAD", "AJAX  ", "", "🔒"),
  @(25.0, "Detector de Metales con contador estadístico de 18 Zonas ", 1139.6153846153845, 8889.0, 7.0, "CONTROL DE ACCESO Y SEGURIDAD", "ZKTECO", "", "🔒"),
  @(26.0, "Energizador para Cerca Eléctrica convencional, YONUSA ", 272.5641025641026, 2126.0, 6.0, "CONTROL DE ACCESO Y SEGURIDAD", "YONUSA", "", "🔒"),
  @(27.0, "Expansor de entrada inalámbrica para alarmas ", 29.74, 231.97199999999998, 4.0, "CONTROL DE ACCESO Y SEGURIDAD", "DAHUA", "", "🔒"),
  @(28.0, "BARRERA MARCA ZKTECO DE ALETA BIOMÉTRICA ", 3307.6923076923076, 25800.0, 6.0, "CONTROL DE ACCESO Y SEGURIDAD", "ZKTECO", "", "🔒"),
  @(29.0, "CERRADURA MARCA ZKTECO PARA PUERTAS DE CRISTAL ", 203.71794871794873, 1589.0, 69.0, "CONTROL DE ACCESO Y SEGURIDAD", "ZKTECO", "", "🔒"),
  @(30.0, "Energizador para Cerca Eléctrica de Alta frecuencia con Interface, YONUSA ", 341.4102564102564, 2663.0, 9.0, "CONTROL DE ACCESO Y SEGURIDAD", "YONUSA", "", "🔒"),
  @(31.0, "Router Dahua de doble antena 300 Mbps ", 22.94871794871795, 179.0, 8.0, "SISTEMAS DE RED", "DAHUA", "", "🌐"),
  @(32.0, "Router IMOU para cámaras Wi-Fi, encriptado WPS, IPv6  ", 17.307692307692307, 135.0, 5.0, "SISTEMAS DE RED", "IMOU", "", "🌐"),
  @(33.0, "Router inalámbrico N a 300 Mbps ", 19.871794871794872, 155.0, 9.0, "SISTEMAS DE RED", "TP-LINK ", "", "🌐"),
  @(34.0, "Router inalámbrico AX1800 con tecnología inalámbrica de sexta generación ", 84.48717948717949, 659.0, 9.0, "SISTEMAS DE RED", "DAHUA", "", "🌐"),
  @(35.0, "DAHUA Switch PoE Ethernet de 26 Puertos ", 284.87179487179486, 2222.0, 9.0, "SISTEMAS DE RED", "DAHUA", "", "🌐"),
  @(36.0, "SWITCH 24 PUERTOS POE FE + 1 PUERTO GIGABIT ", 217.95, 1700.0, 8.0, "SISTEMAS DE RED", "DAHUA", "", "🌐"),
  @(37.0, "Switch Easy Smart de 8 puertos Gigabit con 4 puertos PoE 8 puertos RJ45 10/100/1000Mbps ", 63.46153846153846, 495.0, 9.0, "SISTEMAS DE RED", "TP-LINK ", "", "🌐"),
  @(38.0, "Switch IMOU gigabit duplexor, de 8 puertos, SG108C, RJ45, Puerto MAC, dirección y autoaprendizaje  ", 26.41, 206.0, 2.0, "SISTEMAS DE RED", "IMOU", "", "🌐"),
  @(39.0, "Switch Marca ZKTECO PoE de 24 Puertos 10/100Mbps ", 179.35897435897436, 1399.0, 9.0, "SISTEMAS DE RED", "ZKTECO", "", "🌐"),
  @(40.0, "asdfasdf", 1.2820512820512822, 10.0, 8.0, "SISTEMA DE ALMACENAMIENTO", "as", "", "💾")
)

$rowNum = 2
foreach ($row in $data) {
  $ws.Cells.Item($rowNum, 1).Value = $row[0]
  $ws.Cells.Item($rowNum, 2).Value = $row[1]
  $ws.Cells.Item($rowNum, 3).Value = $row[2]
  $ws.Cells.Item($rowNum, 4).Value = $row[3]
  $ws.Cells.Item($rowNum, 5).Value = $row[4]
  $ws.Cells.Item($rowNum, 6).Value = $row[5]
  $ws.Cells.Item($rowNum, 7).Value = $row[6]
  $ws.Cells.Item($rowNum, 9).Value = $row[8]
  $rowNum = $rowNum + 1
}

# Row 12 (product ID 12) previously held placeholder/test data in column H
# ("gaymer"); the refreshed product row carries no tag, same as every other row.
$ws.Cells.Item(12, 8).ClearContents()

Write-Host "Applied $($data.Count) rows"